# chore: update Sheets via scheduled runner
# Applies updated currentAveragePrice / LevePrice / LeveProfit figures
# (columns H, I, J, K, L, M, N) for a set of leve rows across the
# ALC, ARM, BSM, CRP, CUL, GSM, LTW and WVR sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H45").Value = 4700
$ws.Range("J45").Value = 5250
$ws.Range("L45").Value = 15750
$ws.Range("N45").Value = -16134

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1175.3024
$ws.Range("I61").Value = 1190.6316
$ws.Range("J61").Value = 1058.8
$ws.Range("K61").Value = 1190.6316
$ws.Range("L61").Value = 1058.8
$ws.Range("M61").Value = -978.6315999999999
$ws.Range("N61").Value = -1482.8
$ws.Range("H136").Value = 1175.3024
$ws.Range("I136").Value = 1190.6316
$ws.Range("J136").Value = 1058.8
$ws.Range("K136").Value = 3571.8948
$ws.Range("L136").Value = 3176.4
$ws.Range("M136").Value = -1021.8948
$ws.Range("N136").Value = -8276.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H135").Value = 34250
$ws.Range("J135").Value = 34250
$ws.Range("L135").Value = 34250
$ws.Range("N135").Value = -44390
$ws.Range("H137").Value = 50858
$ws.Range("J137").Value = 50858
$ws.Range("L137").Value = 50858
$ws.Range("N137").Value = -61058
$ws.Range("H138").Value = 36030.668
$ws.Range("J138").Value = 36030.668
$ws.Range("L138").Value = 36030.668
$ws.Range("N138").Value = -46310.668
$ws.Range("H140").Value = 66000
$ws.Range("J140").Value = 66000
$ws.Range("L140").Value = 66000
$ws.Range("N140").Value = -76360

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H32").Value = 538.46155
$ws.Range("I32").Value = 538.46155
$ws.Range("K32").Value = 538.46155
$ws.Range("M32").Value = -222.46155

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H63").Value = 5075
$ws.Range("I63").Value = 3950
$ws.Range("J63").Value = 5450
$ws.Range("K63").Value = 11850
$ws.Range("L63").Value = 16350
$ws.Range("M63").Value = -11101
$ws.Range("N63").Value = -17848
$ws.Range("H64").Value = 6259.2173
$ws.Range("I64").Value = 834.8570999999999
$ws.Range("J64").Value = 8632.375
$ws.Range("K64").Value = 2504.5713
$ws.Range("L64").Value = 25897.125
$ws.Range("M64").Value = -2234.5713
$ws.Range("N64").Value = -26437.125
$ws.Range("H66").Value = 5075
$ws.Range("I66").Value = 3950
$ws.Range("J66").Value = 5450
$ws.Range("K66").Value = 35550
$ws.Range("L66").Value = 49050
$ws.Range("M66").Value = -31806
$ws.Range("N66").Value = -56538
$ws.Range("H67").Value = 6259.2173
$ws.Range("I67").Value = 834.8570999999999
$ws.Range("J67").Value = 8632.375
$ws.Range("K67").Value = 2504.5713
$ws.Range("L67").Value = 25897.125
$ws.Range("M67").Value = -1568.5713
$ws.Range("N67").Value = -27769.125
$ws.Range("H68").Value = 2245.125
$ws.Range("I68").Value = 2245.125
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 6735.375
$ws.Range("L68").Value = 0
$ws.Range("M68").Value = -5924.375
$ws.Range("N68").ClearContents()
$ws.Range("H69").Value = 2457.318
$ws.Range("I69").Value = 487.42856
$ws.Range("J69").Value = 3376.6
$ws.Range("K69").Value = 1462.28568
$ws.Range("L69").Value = 10129.8
$ws.Range("M69").Value = -651.28568
$ws.Range("N69").Value = -11751.8
$ws.Range("H70").Value = 6293
$ws.Range("I70").Value = 4624.8887
$ws.Range("J70").Value = 7794.3
$ws.Range("K70").Value = 13874.6661
$ws.Range("L70").Value = 23382.9
$ws.Range("M70").Value = -13559.6661
$ws.Range("N70").Value = -24012.9
$ws.Range("H71").Value = 2245.125
$ws.Range("I71").Value = 2245.125
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 20206.125
$ws.Range("L71").Value = 0
$ws.Range("M71").Value = -16150.125
$ws.Range("N71").ClearContents()
$ws.Range("H72").Value = 2457.318
$ws.Range("I72").Value = 487.42856
$ws.Range("J72").Value = 3376.6
$ws.Range("K72").Value = 4386.85704
$ws.Range("L72").Value = 30389.4
$ws.Range("M72").Value = -330.8570399999999
$ws.Range("N72").Value = -38501.39999999999
$ws.Range("H73").Value = 6293
$ws.Range("I73").Value = 4624.8887
$ws.Range("J73").Value = 7794.3
$ws.Range("K73").Value = 13874.6661
$ws.Range("L73").Value = 23382.9
$ws.Range("M73").Value = -12782.6661
$ws.Range("N73").Value = -25566.9
$ws.Range("H74").Value = 11565
$ws.Range("I74").Value = 1300
$ws.Range("J74").Value = 14986.667
$ws.Range("K74").Value = 3900
$ws.Range("L74").Value = 44960.001
$ws.Range("M74").Value = -2839
$ws.Range("N74").Value = -47082.001
$ws.Range("H75").Value = 11500
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 11500
$ws.Range("K75").Value = 0
$ws.Range("L75").Value = 34500
$ws.Range("M75").ClearContents()
$ws.Range("N75").Value = -36496
$ws.Range("H76").Value = 7066.6665
$ws.Range("I76").Value = 4000
$ws.Range("K76").Value = 12000
$ws.Range("M76").Value = -11617
$ws.Range("H77").Value = 11565
$ws.Range("I77").Value = 1300
$ws.Range("J77").Value = 14986.667
$ws.Range("K77").Value = 11700
$ws.Range("L77").Value = 134880.003
$ws.Range("M77").Value = -6396
$ws.Range("N77").Value = -145488.003
$ws.Range("H78").Value = 11500
$ws.Range("I78").Value = 0
$ws.Range("J78").Value = 11500
$ws.Range("K78").Value = 0
$ws.Range("L78").Value = 103500
$ws.Range("M78").ClearContents()
$ws.Range("N78").Value = -113484
$ws.Range("H79").Value = 7066.6665
$ws.Range("I79").Value = 4000
$ws.Range("K79").Value = 12000
$ws.Range("M79").Value = -10674
$ws.Range("H80").Value = 7122.222
$ws.Range("J80").Value = 8866.666999999999
$ws.Range("L80").Value = 26600.001
$ws.Range("N80").Value = -28472.001
$ws.Range("H81").Value = 125001384
$ws.Range("J81").Value = 142858640
$ws.Range("L81").Value = 428575920
$ws.Range("N81").Value = -428578166
$ws.Range("H82").Value = 6197.5454
$ws.Range("I82").Value = 1000
$ws.Range("J82").Value = 7352.5557
$ws.Range("K82").Value = 3000
$ws.Range("L82").Value = 22057.6671
$ws.Range("M82").Value = -2594
$ws.Range("N82").Value = -22869.6671
$ws.Range("H83").Value = 7122.222
$ws.Range("J83").Value = 8866.666999999999
$ws.Range("L83").Value = 79800.003
$ws.Range("N83").Value = -89160.003
$ws.Range("H84").Value = 125001384
$ws.Range("J84").Value = 142858640
$ws.Range("L84").Value = 1285727760
$ws.Range("N84").Value = -1285738992
$ws.Range("H85").Value = 6197.5454
$ws.Range("I85").Value = 1000
$ws.Range("J85").Value = 7352.5557
$ws.Range("K85").Value = 3000
$ws.Range("L85").Value = 22057.6671
$ws.Range("M85").Value = -1596
$ws.Range("N85").Value = -24865.6671
$ws.Range("H86").Value = 1500
$ws.Range("I86").Value = 1100
$ws.Range("J86").Value = 1557.1428
$ws.Range("K86").Value = 3300
$ws.Range("L86").Value = 4671.428400000001
$ws.Range("M86").Value = -2114
$ws.Range("N86").Value = -7043.428400000001
$ws.Range("H87").Value = 10103.5
$ws.Range("I87").Value = 7165.6
$ws.Range("K87").Value = 21496.8
$ws.Range("M87").Value = -20248.8
$ws.Range("H88").Value = 4288
$ws.Range("J88").Value = 4288
$ws.Range("L88").Value = 12864
$ws.Range("N88").Value = -13720
$ws.Range("H89").Value = 1500
$ws.Range("I89").Value = 1100
$ws.Range("J89").Value = 1557.1428
$ws.Range("K89").Value = 9900
$ws.Range("L89").Value = 14014.2852
$ws.Range("M89").Value = -3972
$ws.Range("N89").Value = -25870.2852
$ws.Range("H90").Value = 10103.5
$ws.Range("I90").Value = 7165.6
$ws.Range("K90").Value = 64490.4
$ws.Range("M90").Value = -58250.4
$ws.Range("H91").Value = 4288
$ws.Range("J91").Value = 4288
$ws.Range("L91").Value = 12864
$ws.Range("N91").Value = -15828

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H7").Value = 7526000
$ws.Range("J7").Value = 52000
$ws.Range("L7").Value = 52000
$ws.Range("N7").Value = -52224
$ws.Range("H8").Value = 7526000
$ws.Range("J8").Value = 52000
$ws.Range("L8").Value = 52000
$ws.Range("N8").Value = -52278
$ws.Range("H12").Value = 500050
$ws.Range("I12").Value = 100
$ws.Range("J12").Value = 1000000
$ws.Range("K12").Value = 100
$ws.Range("L12").Value = 1000000
$ws.Range("M12").Value = 40
$ws.Range("N12").Value = -1000280
$ws.Range("H51").Value = 14661.5
$ws.Range("J51").Value = 14661.5
$ws.Range("L51").Value = 14661.5
$ws.Range("N51").Value = -15679.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 2434.889
$ws.Range("I122").Value = 1965.2307
$ws.Range("K122").Value = 5895.6921
$ws.Range("M122").Value = -3445.6921

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H109").Value = 22192.334
$ws.Range("J109").Value = 22192.334
$ws.Range("L109").Value = 22192.334
$ws.Range("N109").Value = -24966.334
